$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.074.90"
$ws.Range('E2').Value = '  -2.95%  '
$ws.Range('D3').Value = "'1.842.86"
$ws.Range('E3').Value = '  -2.11%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'0.6991"
$ws.Range('E5').Value = '  -6.08%  '
$ws.Range('D6').Value = "'237.41"
$ws.Range('E6').Value = '  -2.31%  '
$ws.Range('D7').Value = "'1.001"
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = "'0.3032"
$ws.Range('E8').Value = '  -4.12%  '
$ws.Range('D9').Value = "'0.07420"
$ws.Range('E9').Value = '  +2.76%  '
$ws.Range('D10').Value = "'23.22"
$ws.Range('E10').Value = '  -6.80%  '
$ws.Range('D11').Value = "'0.08101"
$ws.Range('E11').Value = '  -2.85%  '
$ws.Range('D12').Value = "'0.7219"
$ws.Range('D13').Value = "'5.228"
$ws.Range('E13').Value = '  -3.51%  '
$ws.Range('D14').Value = "'1.797.48"
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').Value = "'88.89"
$ws.Range('E15').Value = '  -3.94%  '
$ws.Range('D16').Value = "'29.008.78"
$ws.Range('E16').Value = '  -3.29%  '
$ws.Range('D17').Value = "'5.782"
$ws.Range('E17').Value = '  -6.30%  '
$ws.Range('D18').Value = "'240.04"
$ws.Range('E18').Value = '  -3.82%  '
$ws.Range('D19').Value = "'0.000007648"
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('D20').Value = "'12.98"
$ws.Range('E20').Value = '  -4.56%  '
$ws.Range('D21').Value = "'1.0000"
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').Value = "'1.001"
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('B23').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D23').Value = "'2.069.88"
$ws.Range('E23').Value = '  -5.14%  '
$ws.Range('D24').Value = "'7.534"
$ws.Range('E24').Value = '  -6.02%  '
$ws.Range('D25').Value = "'161.27"
$ws.Range('E25').Value = '  -2.53%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = "'8.932"
$ws.Range('E26').Value = '  -3.86%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = "'0.1454"
$ws.Range('E27').Value = '  -6.91%  '
$ws.Range('D28').Value = "'17.95"
$ws.Range('E28').Value = '  -4.07%  '
$ws.Range('D29').Value = "'1.932"
$ws.Range('E29').Value = '  -5.44%  '
$ws.Range('D30').Value = "'1.373"
$ws.Range('E30').Value = '  -7.61%  '
$ws.Range('D31').Value = "'4.472"
$ws.Range('E31').Value = '  -3.15%  '
$ws.Range('E32').Value = '  -3.07%  '
$ws.Range('D33').Value = "'4.011"
$ws.Range('E33').Value = '  -5.27%  '
$ws.Range('D34').Value = "'0.05161"
$ws.Range('E34').Value = '  -3.95%  '
$ws.Range('D35').Value = "'1.181"
$ws.Range('E35').Value = '  -5.86%  '
$ws.Range('D36').Value = "'0.7074"
$ws.Range('E36').Value = '  -6.66%  '
$ws.Range('D37').Value = "'0.9938"
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').Value = "'2.647"
$ws.Range('E38').Value = '  -2.05%  '
$ws.Range('D39').Value = "'0.01866"
$ws.Range('E39').Value = '  -5.17%  '
$ws.Range('D40').Value = "'2.669"
$ws.Range('E40').Value = '  -3.40%  '
$ws.Range('D41').Value = "'0.8956"
$ws.Range('E41').Value = '  +2.90%  '
$ws.Range('D42').Value = "'5.929"
$ws.Range('E42').Value = '  -2.13%  '
$ws.Range('D43').Value = "'0.4271"
$ws.Range('E43').Value = '  -6.33%  '
$ws.Range('D44').Value = "'1.056.70"
$ws.Range('E44').Value = '  -4.67%  '
$ws.Range('D45').Value = "'69.67"
$ws.Range('E45').Value = '  -3.95%  '
$ws.Range('D46').Value = "'0.9998"
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').Value = "'101.34"
$ws.Range('E47').Value = '  -3.11%  '
$ws.Range('D48').Value = "'1.747"
$ws.Range('E48').Value = '  -6.35%  '
$ws.Range('D49').Value = "'1.981.81"
$ws.Range('E49').Value = '  -5.84%  '
$ws.Range('D50').Value = "'9.153"
$ws.Range('E50').Value = '  -4.18%  '
$ws.Range('D51').Value = "'7.016"
$ws.Range('E51').Value = '  -7.85%  '
